$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell I4 (2020), same look as H4 ---
$ws.Cells.Item(4, 8).Copy()
$ws.Cells.Item(4, 9).PasteSpecial(-4122)
$ws.Cells.Item(4, 9).Value = 2020

# --- I5: bold summary row value, like H5 but formatted with one decimal ---
$ws.Cells.Item(5, 8).Copy()
$ws.Cells.Item(5, 9).PasteSpecial(-4122)
$ws.Cells.Item(5, 9).Value = 1.5
$ws.Cells.Item(5, 9).NumberFormat = "0.0"

# --- I6:I13 body rows, like H6:H13 but formatted with one decimal ---
$bodyValues = @{ 6 = 0.2; 7 = 0.8; 8 = 0.4; 9 = 1.8; 10 = 0.5; 11 = 0.7; 12 = 1.9; 13 = 4.6 }
foreach ($r in 6..13) {
  $ws.Cells.Item($r, 8).Copy()
  $ws.Cells.Item($r, 9).PasteSpecial(-4122)
  $ws.Cells.Item($r, 9).Value = $bodyValues[$r]
  $ws.Cells.Item($r, 9).NumberFormat = "0.0"
}

# --- I14: bottom bordered row, formatted with one decimal, no wrap/horizontal align ---
$ws.Cells.Item(14, 8).Copy()
$ws.Cells.Item(14, 9).PasteSpecial(-4122)
$ws.Cells.Item(14, 9).Value = 0.4
$ws.Cells.Item(14, 9).NumberFormat = "0.0"
$ws.Cells.Item(14, 9).HorizontalAlignment = 1
$ws.Cells.Item(14, 9).WrapText = $false

# --- Final selection left on M9, matching the saved view state ---
$ws.Range("M9").Select()

Write-Host "done"
